# ----------------------------------------------------------------------------
# Weekly 84th Precinct CompStat report refresh: new crime data collected.
# Bumps the report volume/number + covered week, and refreshes every stat in
# the "Crime Complaints" table (rows 14-30) with this week's figures.
# ----------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead: volume/issue number and the date range covered by the report ---
$ws.Range("A8").Value = "Volume 30   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/17/2023  Through  7/23/2023"

# --- Row 15 ---
$ws.Range("N15").Value = -80

# --- Row 16 ---
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -61.111111111111
$ws.Range("I16").Value = 72
$ws.Range("J16").Value = 87
$ws.Range("K16").Value = -17.241379310344
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -34.545454545454
$ws.Range("N16").Value = -89.504373177842

# --- Row 17 ---
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 27
$ws.Range("G17").Value = 22
$ws.Range("H17").Value = 22.727272727272
$ws.Range("I17").Value = 142
$ws.Range("J17").Value = 101
$ws.Range("K17").Value = 40.59405940594
$ws.Range("L17").Value = 73.170731707317
$ws.Range("M17").Value = 129.032258064516
$ws.Range("N17").Value = -36.322869955157

# --- Row 18 ---
$ws.Range("C18").Value = 9
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = 125
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 18.181818181818
$ws.Range("I18").Value = 129
$ws.Range("J18").Value = 116
$ws.Range("K18").Value = 11.206896551724
$ws.Range("L18").Value = 98.461538461538
$ws.Range("M18").Value = 98.461538461538
$ws.Range("N18").Value = -71.5859030837

# --- Row 19 ---
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = 7.142857142857
$ws.Range("F19").Value = 58
$ws.Range("G19").Value = 56
$ws.Range("H19").Value = 3.571428571428
$ws.Range("I19").Value = 383
$ws.Range("J19").Value = 351
$ws.Range("K19").Value = 9.116809116809
$ws.Range("L19").Value = 41.328413284132
$ws.Range("M19").Value = 49.027237354085
$ws.Range("N19").Value = -35.304054054054

# --- Row 20 ---
$ws.Range("C20").Value = 1
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 166.666666666667
$ws.Range("I20").Value = 42
$ws.Range("K20").Value = 40
$ws.Range("L20").Value = 147.058823529412
$ws.Range("M20").Value = 50
$ws.Range("N20").Value = -89.447236180904

# --- Row 21 ---
$ws.Range("C21").Value = 32
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 45.454545454545
$ws.Range("F21").Value = 113
$ws.Range("G21").Value = 110
$ws.Range("H21").Value = 2.727272727272
$ws.Range("I21").Value = 771
$ws.Range("J21").Value = 690
$ws.Range("K21").Value = 11.739130434782
$ws.Range("L21").Value = 50
$ws.Range("M21").Value = 47.137404580152
$ws.Range("N21").Value = -67.42712294043

# --- Row 22 ---
$ws.Range("C22").NumberFormat = "General"
$ws.Range("C22").Value = "'0"
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = -100
$ws.Range("G22").Value = 14
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 36
$ws.Range("K22").Value = -22.222222222222
$ws.Range("L22").Value = 3.703703703703
$ws.Range("M22").Value = -26.315789473684

# --- Row 23 ---
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = -20
$ws.Range("J23").Value = 23
$ws.Range("K23").Value = 13.043478260869
$ws.Range("L23").Value = -7.142857142857
$ws.Range("M23").Value = 100

# --- Row 24 ---
$ws.Range("C24").Value = 58
$ws.Range("D24").Value = 48
$ws.Range("E24").Value = 20.833333333333
$ws.Range("F24").Value = 184
$ws.Range("G24").Value = 172
$ws.Range("H24").Value = 6.976744186046
$ws.Range("I24").Value = 1148
$ws.Range("J24").Value = 994
$ws.Range("K24").Value = 15.492957746478
$ws.Range("L24").Value = 51.052631578947
$ws.Range("M24").Value = 38.647342995169

# --- Row 25 ---
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 9
$ws.Range("E25").Value = -44.444444444444
$ws.Range("F25").Value = 33
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = 57.142857142857
$ws.Range("I25").Value = 209
$ws.Range("J25").Value = 202
$ws.Range("K25").Value = 3.465346534653
$ws.Range("L25").Value = 44.137931034482
$ws.Range("M25").Value = 4.5

# --- Row 26 ---
$ws.Range("C26").NumberFormat = "#,##0"
$ws.Range("C26").Value = 1
$ws.Range("F26").Value = 2
$ws.Range("I26").Value = 5
$ws.Range("K26").Value = -37.5
$ws.Range("L26").Value = -61.538461538461

# --- Row 27 ---
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 6
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = -14.285714285714
$ws.Range("I27").Value = 29
$ws.Range("J27").Value = 39
$ws.Range("K27").Value = -25.641025641025
$ws.Range("L27").Value = -12.121212121212

# --- Row 28 ---
$ws.Range("F28").NumberFormat = "General"
$ws.Range("F28").Value = "'0"
$ws.Range("H28").Value = -100

# --- Row 29 ---
$ws.Range("F29").NumberFormat = "General"
$ws.Range("F29").Value = "'0"
$ws.Range("H29").Value = -100

# --- Row 30 ---
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 1
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E30").Value = -100
$ws.Range("J30").Value = 9
$ws.Range("K30").Value = 22.222222222222
